# Fruta / hortaliza, semanal
# A new weekly price record for "Perejil" at "Vega Modelo de Temuco" is
# inserted as row 311 (pushing the existing rows 311:342 down to 312:343).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at 311, shifting rows 311:342 down to 312:343.
$ws.Rows.Item(311).Insert()

# Populate the newly inserted row with the new observation.
$ws.Cells.Item(311, 1).Value  = 10
$ws.Cells.Item(311, 2).Value  = "Vega Modelo de Temuco"
$ws.Cells.Item(311, 3).Value  = "La Araucanía"
$ws.Cells.Item(311, 4).Value  = 44769
$ws.Cells.Item(311, 5).Value  = 9
$ws.Cells.Item(311, 6).Value  = 100112044
$ws.Cells.Item(311, 7).Value  = "Perejil"
$ws.Cells.Item(311, 8).Value  = "Sin especificar"
$ws.Cells.Item(311, 9).Value  = "Primera"
$ws.Cells.Item(311, 10).Value = 50
$ws.Cells.Item(311, 11).Value = 4300
$ws.Cells.Item(311, 12).Value = 4300
$ws.Cells.Item(311, 13).Value = 4300
$ws.Cells.Item(311, 14).Value = "$/docena de atados (3 kilos)"
$ws.Cells.Item(311, 15).Value = "Región Metropolitana"
$ws.Cells.Item(311, 16).Value = 1433
$ws.Cells.Item(311, 17).Value = 3
$ws.Cells.Item(311, 18).Value = "Hortaliza"
